$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.458.63'
$ws.Range("E2").Value = '  -2.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.054.28'
$ws.Range("E3").Value = '  -0.47%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.57'
$ws.Range("E5").Value = '  -3.22%  '

$ws.Range("E6").Value = '  -0.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.24'
$ws.Range("E8").Value = '  -7.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '58.17'
$ws.Range("E9").Value = '  -3.94%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.356'
$ws.Range("E10").Value = '  -7.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0748'
$ws.Range("E11").Value = '  -5.30%  '

$ws.Range("E12").Value = '  -2.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.894'
$ws.Range("E13").Value = '  -2.96%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.59'
$ws.Range("E14").Value = '  -8.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.355.97'
$ws.Range("E15").Value = '  -0.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.34'
$ws.Range("E16").Value = '  -8.90%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.037.73'
$ws.Range("E17").Value = '  -1.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.424.87'
$ws.Range("E18").Value = '  -2.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.60'
$ws.Range("E19").Value = '  -11.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.91'
$ws.Range("E20").Value = '  -4.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0855'
$ws.Range("E21").Value = '  -6.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '237.46'
$ws.Range("E22").Value = '  -0.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.24'
$ws.Range("E23").Value = '  -5.31%  '

$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  -5.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.31'
$ws.Range("E26").Value = '  -3.60%  '

$ws.Range("E27").Value = '  -5.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.51'
$ws.Range("E28").Value = '  -5.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.06'
$ws.Range("E29").Value = '  -1.19%  '

$ws.Range("E30").Value = '  -3.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.03'
$ws.Range("E31").Value = '  -9.33%  '

$ws.Range("E32").Value = '  +0.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.47'
$ws.Range("E33").Value = '  -8.27%  '

$ws.Range("E34").Value = '  -6.36%  '

$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("E36").Value = '  +1.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0827'
$ws.Range("E37").Value = '  -6.63%  '

$ws.Range("E38").Value = '  -7.67%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.24'
$ws.Range("E39").Value = '  -8.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.83'
$ws.Range("E40").Value = '  -7.24%  '

$ws.Range("E41").Value = '  -5.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.83'
$ws.Range("E42").Value = '  -9.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.11'
$ws.Range("E43").Value = '  -5.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.30'
$ws.Range("E44").Value = '  -7.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0896'
$ws.Range("E45").Value = '  -11.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.375.41'
$ws.Range("E46").Value = '  +4.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.64'
$ws.Range("E47").Value = '  -11.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.26'
$ws.Range("E48").Value = '  +3.62%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.85'
$ws.Range("E49").Value = '  -1.21%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.25'
$ws.Range("E50").Value = '  -7.83%  '

$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.242.23'
$ws.Range("E51").Value = '  -0.39%  '
